# Split the leading "{" character off of two single-run field markers so
# that each becomes its own <w:r> (matching the
# TokenIteratorFieldRewriterSplit output): "{m" -> "{" + "m" and
# "{m:" -> "{" + "m:". Re-assigning a Range's FormattedText to itself is a
# content/formatting no-op in Word, but it forces the run containing the
# start of the range to be split at the range boundary without touching
# (or fabricating) any run properties.

$d = $word.ActiveDocument

function Split-RunAtStart($range) {
    $point = $d.Range($range.Start, $range.Start + 1)
    $formatted = $point.FormattedText
    $point.FormattedText = $formatted
}

# Target 1: the "Titre1"/"Heading 1" paragraph whose whole run text is
# "{m" immediately followed by a separate run ":v.name}".
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Heading 1" -and $p.Range.Text -like "{m:v.name}*") {
        Split-RunAtStart $p.Range
    }
}

# Target 2: the paragraph whose run text is "{m:" immediately followed by
# a separate run "endfor}".
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "{m:endfor}*") {
        Split-RunAtStart $p.Range
    }
}
